$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2306
$ws.Range("K3").Value = 2215
$ws.Range("F4").Value = 1908
$ws.Range("I4").Value = 1787
$ws.Range("K5").Value = 147
$ws.Range("K6").Value = 2785
$ws.Range("F7").Value = 24101
$ws.Range("I7").Value = 26241
$ws.Range("K7").Value = 7922

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 157
$ws.Range("K3").Value = 155
$ws.Range("K6").Value = 180
$ws.Range("K7").Value = 532

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 68
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 116
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 65
$ws.Range("K5").Value = 12
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 233
$ws.Range("K8").Value = 532
$ws.Range("K11").Value = 170
$ws.Range("K15").Value = 79
$ws.Range("K19").Value = 229
$ws.Range("K20").Value = 174
$ws.Range("K27").Value = 86
$ws.Range("K29").Value = 399
$ws.Range("K31").Value = 89
$ws.Range("K32").Value = 13
$ws.Range("K33").Value = 311
$ws.Range("K37").Value = 255
$ws.Range("K41").Value = 70
$ws.Range("K42").Value = 272
$ws.Range("K43").Value = 72
$ws.Range("K44").Value = 74
$ws.Range("K49").Value = 53
$ws.Range("K52").Value = 216
$ws.Range("K54").Value = 150
$ws.Range("F63").Value = 193
$ws.Range("I63").Value = 199
$ws.Range("K63").Value = 28
$ws.Range("K64").Value = 52
$ws.Range("K67").Value = 304
$ws.Range("K73").Value = 78
$ws.Range("K76").Value = 114
$ws.Range("K77").Value = 55
$ws.Range("K78").Value = 105
$ws.Range("K79").Value = 208
$ws.Range("K80").Value = 28
$ws.Range("K83").Value = 172
$ws.Range("K85").Value = 388
$ws.Range("K89").Value = 107
$ws.Range("K94").Value = 93
$ws.Range("K95").Value = 125
$ws.Range("K96").Value = 108
$ws.Range("K99").Value = 143
$ws.Range("F101").Value = 24101
$ws.Range("I101").Value = 26241
$ws.Range("K101").Value = 7922

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 32
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 304

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 107
$ws.Range("K7").Value = 399

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 75
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 67
$ws.Range("K7").Value = 272

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 67
$ws.Range("K3").Value = 74
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 50
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K4").Value = 8
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 23
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 141
$ws.Range("K7").Value = 388

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 49
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 216
